# edit.ps1 - apply the certificate-template update:
#  1. Delete the 3rd slide (old "assisted in the coordination..." cover slide)
#  2. Refresh the "Updated: <date>" footer field (slideMaster + all 11 layouts)
#     from 9/8/2017 -> 19/3/2019
#  3. Fix the workshop title wording and dates on the two remaining certificate
#     slides ("GPU programming basics using CUDA" -> "GPU Programming Basics
#     with CUDA", "January 18th to March 1st, 2017" -> "January 15th to
#     March 13th, 2019")

$p = $ppt.ActivePresentation

function Set-DateFieldText($shapes, $oldText, $newText) {
    $updated = 0
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $rng = $shp.TextFrame.TextRange
            if ($rng.Text -eq $oldText) {
                $chars = $rng.Characters(1, $rng.Length)
                $chars.Text = $newText
                $updated = $updated + 1
            }
        }
    }
    return $updated
}

function Update-CertificateSlideText($slideIndex) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(3)
    $range = $shape.TextFrame.TextRange
    $whole = $range.Text

    # -- "GPU programming basics using CUDA " -> "GPU Programming Basics with CUDA "
    $oldTitle = "GPU programming basics using CUDA "
    $newTitle = "GPU Programming Basics with CUDA "
    $titleIdx = $whole.IndexOf($oldTitle)
    if ($titleIdx -ge 0) {
        $chars = $range.Characters($titleIdx + 1, $oldTitle.Length)
        $chars.Text = $newTitle
    }

    # -- re-read text/offsets fresh, then fix the workshop date sentence.
    # The sentence is split across five runs:
    #   "from January 18" + "th" + " to March 1" + "st" + ", 2017."
    $whole = $range.Text
    $anchor = "from January 18"
    $anchorIdx = $whole.IndexOf($anchor)
    if ($anchorIdx -ge 0) {
        $pos = $anchorIdx + 1
        $len1 = $anchor.Length
        $len2 = 2
        $seg3 = " to March 1"
        $len3 = $seg3.Length
        $len4 = 2
        $seg5 = ", 2017."
        $len5 = $seg5.Length

        # apply right-to-left so earlier offsets stay valid while lengths shift
        $c5 = $range.Characters($pos + $len1 + $len2 + $len3 + $len4, $len5)
        $c5.Text = ", 2019."

        $c4 = $range.Characters($pos + $len1 + $len2 + $len3, $len4)
        $c4.Text = "th"

        $c3 = $range.Characters($pos + $len1 + $len2, $len3)
        $c3.Text = " to March 13"

        $c1 = $range.Characters($pos, $len1)
        $c1.Text = "from January 15"
    }
}

# 1. Drop the third (unused) slide entirely.
$p.Slides.Item(3).Delete()

# 2. Fix up the workshop title/date wording on the two certificate slides.
Update-CertificateSlideText 1
Update-CertificateSlideText 2

# 3. Refresh the cached "datetimeFigureOut" footer field everywhere it lives:
#    the slide master and every one of its custom layouts.
$oldDate = "9/8/2017"
$newDate = "19/3/2019"
$master = $p.SlideMaster

$n = Set-DateFieldText $master.Shapes $oldDate $newDate
Write-Output "slideMaster date fields updated: $n"

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    $ln = Set-DateFieldText $layout.Shapes $oldDate $newDate
    Write-Output "layout $L date fields updated: $ln"
}
